$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for "Haba" at Vega Modelo de Temuco.
# It belongs right above the current row 60, so insert a fresh row there
# and let Excel push every row from 60..107 down to 61..108.
$ws.Rows("60:60").Insert()

# Populate the newly inserted row 60 with the new record's data.
$ws.Range("A60").Value = 10
$ws.Range("B60").Value = "Vega Modelo de Temuco"
$ws.Range("C60").Value = "La Araucanía"
$ws.Range("D60").Value = 45264
$ws.Range("E60").Value = 9
$ws.Range("F60").Value = 100112026
$ws.Range("G60").Value = "Haba"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 80
$ws.Range("K60").Value = 13000
$ws.Range("L60").Value = 13000
$ws.Range("M60").Value = 13000
$ws.Range("N60").Value = "`$/saco 25 kilos"
$ws.Range("O60").Value = "Región del Maule"
$ws.Range("P60").Value = 520
$ws.Range("Q60").Value = 25
$ws.Range("R60").Value = "Hortaliza"

# Match the original column's date format (style used by D2:D107).
$ws.Range("D60").NumberFormat = $ws.Range("D61").NumberFormat
